$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 296, shifting existing rows 296..421 down to 297..422.
$ws.Rows("296:296").Insert()

# Copy the constant columns (A,B,C,E,F,G,H,I,J,K,Q,T) from the row that just
# shifted down (297) into the newly-inserted, now-empty row 296. The quality
# column (L) is unchanged by this edit, so it is copied too. Only D, M, N, O,
# P, R, S get brand-new values for this record.
$ws.Range("A296").Value = $ws.Range("A297").Value()
$ws.Range("B296").Value = $ws.Range("B297").Value()
$ws.Range("C296").Value = $ws.Range("C297").Value()
$ws.Range("D296").Value = 44704
$ws.Range("E296").Value = $ws.Range("E297").Value()
$ws.Range("F296").Value = $ws.Range("F297").Value()
$ws.Range("G296").Value = $ws.Range("G297").Value()
$ws.Range("H296").Value = $ws.Range("H297").Value()
$ws.Range("I296").Value = $ws.Range("I297").Value()
$ws.Range("J296").Value = $ws.Range("J297").Value()
$ws.Range("K296").Value = $ws.Range("K297").Value()
$ws.Range("L296").Value = $ws.Range("L297").Value()
$ws.Range("M296").Value = 250
$ws.Range("N296").Value = 7000
$ws.Range("O296").Value = 7500
$ws.Range("P296").Value = 7300
$ws.Range("Q296").Value = $ws.Range("Q297").Value()
$ws.Range("R296").Value = "Brasil"
$ws.Range("S296").Value = 1825
$ws.Range("T296").Value = $ws.Range("T297").Value()
